# Apply updated "想去人数" (F) / "最低票价" (G) figures across the
# workbook's sheets, matching the regenerated gh-pages data snapshot.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsLocal   = $wb.Worksheets.Item("本地生活")
$wsAll     = $wb.Worksheets.Item("全部类型")

# ---- 展览 (Exhibition) sheet ----
$wsExhibit.Range("F2").Value  = 9120
$wsExhibit.Range("F3").Value  = 1983
$wsExhibit.Range("F4").Value  = 6645
$wsExhibit.Range("G6").Value  = 85
$wsExhibit.Range("G7").Value  = 70
$wsExhibit.Range("F15").Value = 38
$wsExhibit.Range("F16").Value = 9169
$wsExhibit.Range("F17").Value = 170
$wsExhibit.Range("F21").Value = 1868
$wsExhibit.Range("F25").Value = 116
$wsExhibit.Range("F29").Value = 33
$wsExhibit.Range("F31").Value = 579
$wsExhibit.Range("F33").Value = 66
$wsExhibit.Range("F35").Value = 2464
$wsExhibit.Range("G35").Value = 69
$wsExhibit.Range("F36").Value = 891
$wsExhibit.Range("F37").Value = 571
$wsExhibit.Range("F41").Value = 336
$wsExhibit.Range("F42").Value = 191
$wsExhibit.Range("F47").Value = 35
$wsExhibit.Range("F48").Value = 4010
$wsExhibit.Range("F49").Value = 19

# ---- 本地生活 (Local life) sheet ----
$wsLocal.Range("F5").Value = 29

# ---- 全部类型 (All types) sheet ----
$wsAll.Range("F3").Value  = 9120
$wsAll.Range("F5").Value  = 1983
$wsAll.Range("F6").Value  = 6645
$wsAll.Range("G8").Value  = 70
$wsAll.Range("F14").Value = 29
$wsAll.Range("F17").Value = 9169
$wsAll.Range("F18").Value = 170
$wsAll.Range("F22").Value = 1868
$wsAll.Range("F24").Value = 116
$wsAll.Range("F27").Value = 33
$wsAll.Range("F29").Value = 579
$wsAll.Range("F31").Value = 66
$wsAll.Range("F33").Value = 891
$wsAll.Range("F36").Value = 571
$wsAll.Range("F37").Value = 336
$wsAll.Range("F39").Value = 191
$wsAll.Range("F44").Value = 35
$wsAll.Range("F45").Value = 4010
$wsAll.Range("F48").Value = 19
